# Rolling-average style shift: each row's Confirmed (H) value and Date (I) label
# are moved back to reflect the PREVIOUS day's figures for that site, the series
# effectively gaining one earlier day (2020-03-18) "for free" and one new day
# (2020-05-07) appended at the end carrying forward the last known counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 99

# Snapshot the current H (Confirmed) and I (Date) columns before mutating anything,
# since writes below depend on the ORIGINAL values of earlier rows.
$origH = @{}
$origI = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $origH[$r] = $ws.Cells.Item($r, 8).Value2
    $origI[$r] = $ws.Cells.Item($r, 9).Value2
}

# New first date = one day before the current earliest date.
$firstDate = [datetime]::ParseExact($origI[$firstRow], "yyyy-MM-dd", $null)
$newFirstDateText = $firstDate.AddDays(-1).ToString("yyyy-MM-dd")

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $prevRow = $r - 2
    if ($prevRow -ge $firstRow) {
        $newH = $origH[$prevRow]
        $newDateText = $origI[$prevRow]
    } else {
        $newH = 0
        $newDateText = $newFirstDateText
    }

    $ws.Cells.Item($r, 8).Value = $newH

    $cellI = $ws.Cells.Item($r, 9)
    $cellI.NumberFormat = "@"
    $cellI.Value = $newDateText
    $cellI.ClearFormats()
}

# Append two new rows for the new final date (one day after the previous
# final date), carrying forward the last (pre-shift) confirmed counts for
# each site.
$oldLastDate = [datetime]::ParseExact($origI[$lastRow], "yyyy-MM-dd", $null)
$newLastDateText = $oldLastDate.AddDays(1).ToString("yyyy-MM-dd")
$kewauneeFinalH = $origH[$lastRow - 1]
$manitowocFinalH = $origH[$lastRow]

$newRows = @(
    @{ Row = 100; Idx = 98; FIPS = 55061; County = "Kewaunee";  Combined = "Kewaunee, Wisconsin, US";  H = $kewauneeFinalH },
    @{ Row = 101; Idx = 99; FIPS = 55071; County = "Manitowoc"; Combined = "Manitowoc, Wisconsin, US"; H = $manitowocFinalH }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    $ws.Range("A" + ($lastRow) + ":A" + ($lastRow)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $nr.Idx

    $ws.Cells.Item($r, 2).Value = $nr.FIPS
    $ws.Cells.Item($r, 3).Value = "Point Beach"
    $ws.Cells.Item($r, 4).Value = 3
    $ws.Cells.Item($r, 5).Value = $nr.County
    $ws.Cells.Item($r, 6).Value = "Wisconsin"
    $ws.Cells.Item($r, 7).Value = $nr.Combined
    $ws.Cells.Item($r, 8).Value = $nr.H

    $cellI = $ws.Cells.Item($r, 9)
    $cellI.NumberFormat = "@"
    $cellI.Value = $newLastDateText
    $cellI.ClearFormats()
}
